$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.055.75'
$ws.Range("E2").Value = '  -2.40%  '
$ws.Range("D3").Value = '3.206.24'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.27'
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.10'
$ws.Range("E6").Value = '  -3.05%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '3.203.90'
$ws.Range("E8").Value = '  -0.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  -3.50%  '
$ws.Range("E10").Value = '  -4.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.59'
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.481'
$ws.Range("E12").Value = '  -5.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000261'
$ws.Range("E13").Value = '  -5.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.40'
$ws.Range("E14").Value = '  -4.31%  '
$ws.Range("D15").Value = '3.726.26'
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").Value = '65.122.68'
$ws.Range("E16").Value = '  -2.44%  '
$ws.Range("D17").Value = '3.197.17'
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.10'
$ws.Range("E19").Value = '  -5.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '486.38'
$ws.Range("E20").Value = '  -5.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.94'
$ws.Range("E21").Value = '  -2.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.725'
$ws.Range("E22").Value = '  -1.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.81'
$ws.Range("E23").Value = '  -3.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.02'
$ws.Range("E24").Value = '  -5.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.40'
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.95'
$ws.Range("E27").Value = '  -1.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.77'
$ws.Range("E28").Value = '  -5.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.129'
$ws.Range("E29").Value = '  +34.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.30'
$ws.Range("E30").Value = '  -4.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.96'
$ws.Range("E31").Value = '  -0.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.76'
$ws.Range("E32").Value = '  -8.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.16'
$ws.Range("E33").Value = '  -3.99%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  -6.15%  '
$ws.Range("E36").Value = '  -5.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.30'
$ws.Range("E37").Value = '  +6.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '54.82'
$ws.Range("E38").Value = '  -2.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '479.80'
$ws.Range("E39").Value = '  -7.00%  '
$ws.Range("D40").Value = '0.0₃0736'
$ws.Range("E40").Value = '  -4.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0408'
$ws.Range("E41").Value = '  -3.33%  '
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.47'
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D45").Value = '2.922.19'
$ws.Range("E45").Value = '  +1.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.281'
$ws.Range("E46").Value = '  -6.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '27.72'
$ws.Range("E47").Value = '  -3.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.38'
$ws.Range("E48").Value = '  -1.60%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.116'
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.01'
$ws.Range("E51").Value = '  -2.03%  '
